$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Con SmartScore"
$ws.Range("B3").Value = "Sin SmartScore"
$ws.Range("B4").Value = "Sin SmartScore"
$ws.Range("B7").Value = "Sin SmartScore"
$ws.Range("B10").Value = "Con SmartScore"
$ws.Range("B15").Value = "Con SmartScore"

$ws.Range("I15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("O15").Value = 0
$ws.Range("R15").Value = 0
$ws.Range("U15").Value = 0
$ws.Range("X15").Value = 0
$ws.Range("AA15").Value = 0
$ws.Range("AD15").Value = 0
$ws.Range("AG15").Value = 0
